$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.215597054981016
$ws.Range("D2").Value = 0.2370224868694484

# Row 3
$ws.Range("C3").Value = 1.801707070589514
$ws.Range("D3").Value = 0.08531326570668618
$ws.Range("G3").Value = "No"

# Row 4
$ws.Range("C4").Value = 1.576522444534856
$ws.Range("D4").Value = 0.129177452807973
$ws.Range("G4").Value = "No"

# Row 5
$ws.Range("C5").Value = 5.265548102927215
$ws.Range("D5").Value = 0.0000277630890566094

# Row 6
$ws.Range("C6").Value = 0.5855778351062695
$ws.Range("D6").Value = 0.5641185150759507

# Row 7
$ws.Range("C7").Value = 0.7349359267559981
$ws.Range("D7").Value = 0.4701408986125641

# Row 8
$ws.Range("C8").Value = 3.126536146024777
$ws.Range("D8").Value = 0.004910503644846465

# Row 9
$ws.Range("C9").Value = -0.1905783498554184
$ws.Range("D9").Value = 0.8506036388228693

# Row 10
$ws.Range("C10").Value = 1.770872366068002
$ws.Range("D10").Value = 0.0904377783540018
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = 2.295389372016389
$ws.Range("D11").Value = 0.03161525401363829
